$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.429.16"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").Value = "2.542.18"
$ws.Range("E3").Value = "  +4.57%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.73"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.47"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +8.12%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").Value = "2.539.75"
$ws.Range("E9").Value = "  +4.68%  "
$ws.Range("E10").Value = "  +1.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.68"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.356"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.19"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +7.56%  "
$ws.Range("D15").Value = "2.996.22"
$ws.Range("E15").Value = "  +4.76%  "
$ws.Range("D16").Value = "63.363.63"
$ws.Range("E16").Value = "  +1.50%  "
$ws.Range("E17").Value = "  +1.32%  "
$ws.Range("D18").Value = "2.620.83"
$ws.Range("E18").Value = "  +7.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.58"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "338.85"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.32"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.77"
$ws.Range("D22").ClearFormats()
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.89"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.83%  "
$ws.Range("E25").Value = "  -1.79%  "
$ws.Range("E26").Value = "  +14.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.58"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.16%  "
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.42"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +3.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.13"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +10.31%  "
$ws.Range("E31").Value = "  +3.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.85"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "177.77"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.98%  "
$ws.Range("E34").Value = "  +9.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "418.95"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +10.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.405"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.99%  "
$ws.Range("E37").Value = "  +2.27%  "
$ws.Range("E38").Value = "  -0.89%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("E40").Value = "  +3.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.46"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "153.37"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +5.90%  "
$ws.Range("E44").Value = "  +2.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.69"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.607"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0962"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0521"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("E49").Value = "  +6.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.42"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.80"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +4.46%  "
